$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "3046010569"
$ws.Range("C4").Value = "732111193280551"
$ws.Range("C8").Select()
